# Applies the three "Sacrifice Summoned" effect-text edits described in the
# commit diff.
#
# Note: Find.Execute's Replace argument runs AutoCorrect/"smart quotes" on
# the replacement text (straight "'" becomes a curly "'"), even with the
# AutoFormat options turned off. To keep the original straight apostrophes
# exactly as in the diff, we locate the target text with Find (no replace)
# and then assign the new text directly via Range.Text, which performs a
# literal, un-autocorrected replacement.

$d = $word.ActiveDocument

# NOTE on ordering: edit 1's replacement text is identical to edit 3's
# *search* text ("Sacrifice Summoned: None of your opponent's cards can be
# activated for 2 of their turns."). Edit 3 must therefore run first,
# while that text still uniquely identifies the original split-run
# paragraph; otherwise Find would match the freshly-written text from
# edit 1 instead of (or in addition to) the real target.

# ---------------------------------------------------------------------
# Edit 3 (applied first): "Sacrifice Summoned: None of your opponent's
#          cards can be activated for 2 of their turns." (spread across
#          4 runs)
#      -> "Sacrifice Summoned: Summon the sacrificed cards. You cannot
#          attack this turn."
# ---------------------------------------------------------------------
$r3 = $d.Content
$found3 = $r3.Find.Execute(
    "Sacrifice Summoned: None of your opponent's cards can be activated for 2 of their turns.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found3) {
    $r3.Text = "Sacrifice Summoned: Summon the sacrificed cards. You cannot attack this turn."
} else {
    Write-Output "WARNING: edit 3 target text not found"
}

# ---------------------------------------------------------------------
# Edit 1: "Bring back the card sacrificed. You cannot attack this turn."
#      -> "None of your opponent's cards can be activated for 2 of their
#          turns."
# ---------------------------------------------------------------------
$r1 = $d.Content
$found1 = $r1.Find.Execute(
    "Sacrifice Summoned: Bring back the card sacrificed. You cannot attack this turn.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $r1.Text = "Sacrifice Summoned: None of your opponent's cards can be activated for 2 of their turns."
} else {
    Write-Output "WARNING: edit 1 target text not found"
}

# ---------------------------------------------------------------------
# Edit 2: append a new run containing "." right after the existing
# "...while this card remains on your battlefield)" run (kept as a
# separate <w:r> rather than merged text).
# ---------------------------------------------------------------------
$r2 = $d.Content
$found2 = $r2.Find.Execute(
    "Sacrifice Summoned: All cards in your battlefield / graveyard / hand cannot be targeted or destroyed for 3 turns (while this card remains on your battlefield)",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $r2.Collapse(0)            # wdCollapseEnd
    $r2.InsertAfter(".")
    # Toggling a character property on the newly-inserted run forces it to
    # stay a distinct <w:r> instead of being coalesced back into the
    # identically-formatted run before it.
    $r2.Bold = 1
    $r2.Bold = 0
} else {
    Write-Output "WARNING: edit 2 target text not found"
}

Write-Output "Done: edit1=$found1 edit2=$found2 edit3=$found3"
